$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: add "Not skinnable" status, replace comment with the fuller explanation.
$ws.Cells.Item(26, 2).Value = "Not skinnable"
$ws.Cells.Item(26, 3).Value = "There is very limited scope for skins to do this (existing context menus ca be added to, but new ones can't be created)"

# Update the view: scrolled position and active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select()
